# Add/update metadata report for Akurana
# Appends a new row (row 3) to the "Metadata Report" sheet for 2025 JAN data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata Report")

$rowNum = 3

$values = @(
    2025,
    "JAN",
    "10/01-01/01",
    "11/01, 12/01, 13/01, 14/01, 15/01, 16/01, 17/01, 18/01, 19/01, 20/01, 21/01, 22/01, 23/01, 24/01, 25/01, 26/01, 27/01, 28/01, 29/01, 30/01, 31/01",
    "Akurana",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "-",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "-",
    "-",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "✓",
    "-",
    "-",
    "-"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rowNum, $col).Value = $values[$i]
}
